# Mifos Automation Excels - "Product mix 5 test cases, Data tables 15 Test Cases"
#
# The underlying automation run cleared out the stale "Linked Savings
# Account (Dividend Posting)" share-id value on the "Share Other Details1"
# sheet (it's populated dynamically by the test each run, so the checked-in
# workbook should not carry a stale id), and left the workbook positioned
# with "Share Other Details2" as the active/selected sheet instead of
# "NewSavingInput".

$wb = $excel.ActiveWorkbook

# --- "Share Other Details1" sheet -----------------------------------------
$wsDetails1 = $wb.Worksheets.Item("Share Other Details1")

# Clear the stale "Linked Savings Account(Dividend Posting)" id value - the
# row shrinks back down to its un-populated wrapped-label height.
$wsDetails1.Range("B5").ClearContents()
$wsDetails1.Rows.Item(5).RowHeight = 45

# Leave the cursor sitting on B10 on this sheet (matches the last recorded
# selection for this tab).
$wsDetails1.Range("B10").Select()

# --- "Share Other Details2" sheet ------------------------------------------
# This becomes the active tab / selected sheet of the workbook, with B2
# selected, replacing "NewSavingInput" as the previously active tab.
$wsDetails2 = $wb.Worksheets.Item("Share Other Details2")
$wsDetails2.Activate()
$wsDetails2.Range("B2").Select()
